# Refresh the cryptocurrency snapshot table (columns B:E, rows 2-51) with
# the latest scraped values, as produced by the scheduled GitHub Actions
# update job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is one cell that needs its text re-written: Coin name (B),
# Link (C), Price (D) and Volume/1h change (E). A couple of rows (15/16
# and 35/36) also swap their coin identity because the ranking order
# changed between scrapes.
$updates = @(
    @{Cell='D2'; Value='20.528.19'; Numeric=$false},
    @{Cell='E2'; Value='  +1.67%  '; Numeric=$false},
    @{Cell='D3'; Value='1.466.93'; Numeric=$false},
    @{Cell='E3'; Value='  +1.95%  '; Numeric=$false},
    @{Cell='D4'; Value='1.007'; Numeric=$true},
    @{Cell='E4'; Value='  +0.20%  '; Numeric=$false},
    @{Cell='D5'; Value='0.9582'; Numeric=$true},
    @{Cell='E5'; Value='  +3.56%  '; Numeric=$false},
    @{Cell='D6'; Value='277.59'; Numeric=$true},
    @{Cell='E6'; Value='  +0.11%  '; Numeric=$false},
    @{Cell='D7'; Value='0.3609'; Numeric=$true},
    @{Cell='E7'; Value='  -0.80%  '; Numeric=$false},
    @{Cell='D8'; Value='0.3077'; Numeric=$true},
    @{Cell='E8'; Value='  -0.43%  '; Numeric=$false},
    @{Cell='D9'; Value='39.37'; Numeric=$true},
    @{Cell='E9'; Value='  +1.11%  '; Numeric=$false},
    @{Cell='D10'; Value='1.072'; Numeric=$true},
    @{Cell='E10'; Value='  +4.50%  '; Numeric=$false},
    @{Cell='D11'; Value='0.06626'; Numeric=$true},
    @{Cell='E11'; Value='  +1.86%  '; Numeric=$false},
    @{Cell='D12'; Value='1.002'; Numeric=$true},
    @{Cell='E12'; Value='  +0.19%  '; Numeric=$false},
    @{Cell='D13'; Value='5.481'; Numeric=$true},
    @{Cell='E13'; Value='  +2.12%  '; Numeric=$false},
    @{Cell='D14'; Value='18.13'; Numeric=$true},
    @{Cell='E14'; Value='  +3.25%  '; Numeric=$false},
    @{Cell='B15'; Value='Chainlink'; Numeric=$false},
    @{Cell='C15'; Value='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Numeric=$false},
    @{Cell='D15'; Value='6.139'; Numeric=$true},
    @{Cell='E15'; Value='  +1.51%  '; Numeric=$false},
    @{Cell='B16'; Value='Dai'; Numeric=$false},
    @{Cell='C16'; Value='https://coinranking.com/coin/MoTuySvg7+dai-dai'; Numeric=$false},
    @{Cell='D16'; Value='0.9580'; Numeric=$true},
    @{Cell='E16'; Value='  +1.70%  '; Numeric=$false},
    @{Cell='D17'; Value='0.00001020'; Numeric=$true},
    @{Cell='E17'; Value='  +0.89%  '; Numeric=$false},
    @{Cell='D18'; Value='1.467.20'; Numeric=$false},
    @{Cell='E18'; Value='  +2.15%  '; Numeric=$false},
    @{Cell='D19'; Value='0.05948'; Numeric=$true},
    @{Cell='E19'; Value='  +5.74%  '; Numeric=$false},
    @{Cell='D20'; Value='68.32'; Numeric=$true},
    @{Cell='E20'; Value='  +0.62%  '; Numeric=$false},
    @{Cell='D21'; Value='5.467'; Numeric=$true},
    @{Cell='E21'; Value='  +1.55%  '; Numeric=$false},
    @{Cell='D22'; Value='14.51'; Numeric=$true},
    @{Cell='E22'; Value='  +1.09%  '; Numeric=$false},
    @{Cell='D23'; Value='11.11'; Numeric=$true},
    @{Cell='E23'; Value='  +2.34%  '; Numeric=$false},
    @{Cell='D24'; Value='2.266'; Numeric=$true},
    @{Cell='E24'; Value='  +1.23%  '; Numeric=$false},
    @{Cell='D25'; Value='20.534.93'; Numeric=$false},
    @{Cell='E25'; Value='  +1.48%  '; Numeric=$false},
    @{Cell='D26'; Value='142.71'; Numeric=$true},
    @{Cell='E26'; Value='  +3.82%  '; Numeric=$false},
    @{Cell='D27'; Value='2.105'; Numeric=$true},
    @{Cell='E27'; Value='  -1.95%  '; Numeric=$false},
    @{Cell='D28'; Value='17.12'; Numeric=$true},
    @{Cell='E28'; Value='  +1.25%  '; Numeric=$false},
    @{Cell='D29'; Value='1.626.54'; Numeric=$false},
    @{Cell='E29'; Value='  +2.38%  '; Numeric=$false},
    @{Cell='D30'; Value='113.58'; Numeric=$true},
    @{Cell='E30'; Value='  +3.11%  '; Numeric=$false},
    @{Cell='D31'; Value='3.893'; Numeric=$true},
    @{Cell='E31'; Value='  +0.66%  '; Numeric=$false},
    @{Cell='D32'; Value='0.07996'; Numeric=$true},
    @{Cell='E32'; Value='  +4.32%  '; Numeric=$false},
    @{Cell='D33'; Value='4.925'; Numeric=$true},
    @{Cell='E33'; Value='  +1.53%  '; Numeric=$false},
    @{Cell='D34'; Value='0.8016'; Numeric=$true},
    @{Cell='E34'; Value='  -1.71%  '; Numeric=$false},
    @{Cell='B35'; Value='TrustWalletToken'; Numeric=$false},
    @{Cell='C35'; Value='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Numeric=$false},
    @{Cell='D35'; Value='1.219'; Numeric=$true},
    @{Cell='E35'; Value='  +6.94%  '; Numeric=$false},
    @{Cell='B36'; Value='WEMIXTOKEN'; Numeric=$false},
    @{Cell='C36'; Value='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; Numeric=$false},
    @{Cell='D36'; Value='1.477'; Numeric=$true},
    @{Cell='E36'; Value='  +0.79%  '; Numeric=$false},
    @{Cell='D37'; Value='0.05760'; Numeric=$true},
    @{Cell='E37'; Value='  -2.36%  '; Numeric=$false},
    @{Cell='D38'; Value='4.688'; Numeric=$true},
    @{Cell='E38'; Value='  +1.16%  '; Numeric=$false},
    @{Cell='D39'; Value='0.02047'; Numeric=$true},
    @{Cell='E39'; Value='  +3.07%  '; Numeric=$false},
    @{Cell='D40'; Value='0.9588'; Numeric=$true},
    @{Cell='E40'; Value='  +3.53%  '; Numeric=$false},
    @{Cell='D41'; Value='10.35'; Numeric=$true},
    @{Cell='E41'; Value='  +1.47%  '; Numeric=$false},
    @{Cell='D42'; Value='0.1865'; Numeric=$true},
    @{Cell='E42'; Value='  +1.37%  '; Numeric=$false},
    @{Cell='D43'; Value='7.372'; Numeric=$true},
    @{Cell='E43'; Value='  +2.34%  '; Numeric=$false},
    @{Cell='D44'; Value='0.5254'; Numeric=$true},
    @{Cell='E44'; Value='  +0.70%  '; Numeric=$false},
    @{Cell='D45'; Value='3.517'; Numeric=$true},
    @{Cell='E45'; Value='  +0.50%  '; Numeric=$false},
    @{Cell='D46'; Value='12.09'; Numeric=$true},
    @{Cell='E46'; Value='  +1.15%  '; Numeric=$false},
    @{Cell='D47'; Value='118.61'; Numeric=$true},
    @{Cell='E47'; Value='  +1.47%  '; Numeric=$false},
    @{Cell='D48'; Value='0.5188'; Numeric=$true},
    @{Cell='E48'; Value='  +1.72%  '; Numeric=$false},
    @{Cell='D49'; Value='1.807'; Numeric=$true},
    @{Cell='E49'; Value='  +3.18%  '; Numeric=$false},
    @{Cell='D50'; Value='0.06438'; Numeric=$true},
    @{Cell='E50'; Value='  +1.83%  '; Numeric=$false},
    @{Cell='D51'; Value='0.9844'; Numeric=$true},
    @{Cell='E51'; Value='  +0.45%  '; Numeric=$false}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # These price strings look like plain numbers (e.g. "1.007",
        # "277.59"), but the sheet stores them as literal text. Assigning
        # a numeric-looking string straight to .Value would make Excel
        # silently reinterpret it as a real number (and round-trip it
        # through floating point), so the cell is forced to Text format
        # first, then the temporary format override is cleared again
        # once the text value is safely in place.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
